# Update the speaker notes on slide 11 ("Reflection") from a single
# "Command and argument, " line to two separate lines answering the
# "easiest" / "most difficult" prompts on that slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$np = $s.NotesPage

# Shape 2 on a notes page is the "Notes Placeholder" body text shape.
$notesShape = $np.Shapes.Item(2)

$notesShape.TextFrame.TextRange.Text = "Easiest: moving and changing directories, `nDifficult: Altering text using nano, exiting different programs, "
